$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 91.1369967353351
$ws.Range("C2").Value = 127.3679339628222
$ws.Range("D2").Value = 145.3129274435881
$ws.Range("E2").Value = 155.734922053209

$ws.Range("B3").Value = 110.2596756929998
$ws.Range("C3").Value = 152.6481160815642
$ws.Range("D3").Value = 171.7948606769434
$ws.Range("E3").Value = 185.1000274471107

$ws.Range("B4").Value = 98.10551322810247
$ws.Range("C4").Value = 141.10191723647
$ws.Range("D4").Value = 162.7965677379082
$ws.Range("E4").Value = 176.7807258803518

$ws.Range("B5").Value = 76.26618851682456
$ws.Range("C5").Value = 103.7042077698815
$ws.Range("D5").Value = 113.4099841640455
$ws.Range("E5").Value = 123.0184044370315

$ws.Range("B6").Value = 67.06515550277362
$ws.Range("C6").Value = 90.96069586813786
$ws.Range("D6").Value = 100.4437278407853
$ws.Range("E6").Value = 107.6300106607524

$ws.Range("B7").Value = 7.210567491758278
$ws.Range("C7").Value = 9.662075413724656
$ws.Range("D7").Value = 10.69148570057467
$ws.Range("E7").Value = 11.397879639144

$ws.Range("B8").Value = 342.9336348477351
$ws.Range("C8").Value = 483.9432969501058
$ws.Range("D8").Value = 563.38277808048
$ws.Range("E8").Value = 608.0822332223942

$ws.Range("B9").Value = 97.93106153690282
$ws.Range("C9").Value = 132.3281512690647
$ws.Range("D9").Value = 146.2917903633075
$ws.Range("E9").Value = 156.9739476974102

$ws.Range("B10").Value = 42.18158432580129
$ws.Range("C10").Value = 55.48839697489341
$ws.Range("D10").Value = 61.96800176324722
$ws.Range("E10").Value = 65.25496252721278

$ws.Range("B11").Value = 7.750659366753535
$ws.Range("C11").Value = 9.681968837488911
$ws.Range("D11").Value = 10.69186417006868
$ws.Range("E11").Value = 12.14826695821297

$ws.Range("B12").Value = 18.5224724752061
$ws.Range("C12").Value = 25.42504606204698
$ws.Range("D12").Value = 28.24126807874219
$ws.Range("E12").Value = 29.35038238180431

$ws.Range("B13").Value = 23.70708953324266
$ws.Range("C13").Value = 30.91395792474116
$ws.Range("D13").Value = 34.85533186473179
$ws.Range("E13").Value = 37.13588811466538

